# Update course excel files
# DIANA SCHOOL OF COMMUNITY SERVICES workbook: split the bogus single
# "DIANA SCHOOL OF COMMUNITY SERVICES" department label into the real
# per-course department names, split the "NSW/QLD/TAS (Currently not
# accepting enrolments)" location into location + locationDetail, and
# drop the (stale) promotion-validity note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- department (column C) ----------------------------------------
$ws.Range("C2:C7").Value  = "Ageing Support"
$ws.Range("C8").Value     = "Community Services"
$ws.Range("C9:C10").Value = "Early Childhood"
$ws.Range("C11:C12").Value = "Massage"
$ws.Range("C13:C20").Value = "Packages"

# ---- location / locationDetail (columns M / N) ---------------------
# Rows 11, 12 and 20 used to cram the "currently not accepting
# enrolments" caveat into the location string itself; split it into
# its own locationDetail cell instead.
$ws.Range("M11").Value = "NSW/QLD/TAS"
$ws.Range("N11").Value = "Currently not accepting enrolments"

$ws.Range("M12").Value = "NSW/QLD/TAS"
$ws.Range("N12").Value = "Currently not accepting enrolments"

$ws.Range("M20").Value = "NSW/QLD/TAS"
$ws.Range("N20").Value = "Currently not accepting enrolments"

# ---- promotionValidity (column R) -----------------------------------
# The "Promotion valid until 31th Dec 2021" note is stale; clear the
# cell contents (formatting/style is left untouched).
$ws.Range("R2:R20").ClearContents()

# ---- row heights -----------------------------------------------------
for ($r = 2; $r -le 17; $r++) {
    $ws.Rows.Item($r).RowHeight = 42.75
}
$ws.Rows.Item(18).RowHeight = 57
$ws.Rows.Item(19).RowHeight = 57
$ws.Rows.Item(20).RowHeight = 42.75

# ---- selection ---------------------------------------------------
$ws.Range("R2:R20").Select()
